$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H4").Value = 136.28572
$ws.Range("J4").Value = 218.66667
$ws.Range("L4").Value = 218.66667
$ws.Range("N4").Value = -446.66667

$ws.Range("H10").Value = 6969
$ws.Range("J10").Value = 6969
$ws.Range("L10").Value = 6969
$ws.Range("N10").Value = -7555

$ws.Range("H13").Value = 1900
$ws.Range("J13").Value = 1900
$ws.Range("L13").Value = 1900
$ws.Range("N13").Value = -2238

$ws.Range("H18").Value = 1900.5
$ws.Range("I18").Value = 301
$ws.Range("J18").Value = 3500
$ws.Range("K18").Value = 301
$ws.Range("L18").Value = 3500
$ws.Range("M18").Value = -17
$ws.Range("N18").Value = -4068

$ws.Range("H32").Value = 999
$ws.Range("I32").Value = 999
$ws.Range("K32").Value = 999
$ws.Range("M32").Value = -673

$ws.Range("H43").Value = 9000.666999999999
$ws.Range("J43").Value = 9000.666999999999
$ws.Range("L43").Value = 9000.666999999999
$ws.Range("N43").Value = -9138.666999999999

$ws.Range("H51").Value = 4500
$ws.Range("I51").Value = 4500
$ws.Range("K51").Value = 4500
$ws.Range("M51").Value = -4016

$ws.Range("H58").Value = 206
$ws.Range("J58").Value = 397
$ws.Range("L58").Value = 1191
$ws.Range("N58").Value = -1491

$ws.Range("H62").Value = 6000
$ws.Range("I62").Value = 0
$ws.Range("K62").Value = 0
$ws.Range("M62").ClearContents()

$ws.Range("H65").Value = 6000
$ws.Range("I65").Value = 0
$ws.Range("K65").Value = 0
$ws.Range("M65").ClearContents()

$ws.Range("H112").Value = 3872.5
$ws.Range("J112").Value = 3872.5
$ws.Range("L112").Value = 11617.5
$ws.Range("N112").Value = -13833.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 96.333336
$ws.Range("I5").Value = 87.833336
$ws.Range("K5").Value = 87.833336
$ws.Range("M5").Value = 24.166664

$ws.Range("H32").Value = 4284.25
$ws.Range("I32").Value = 4284.25
$ws.Range("K32").Value = 4284.25
$ws.Range("M32").Value = -3997.25

$ws.Range("H61").Value = 2000
$ws.Range("I61").Value = 2000
$ws.Range("K61").Value = 2000
$ws.Range("M61").Value = -1788

$ws.Range("H74").Value = 6737.5
$ws.Range("I74").Value = 1975
$ws.Range("K74").Value = 1975
$ws.Range("M74").Value = -1101

$ws.Range("H77").Value = 6737.5
$ws.Range("I77").Value = 1975
$ws.Range("K77").Value = 9875
$ws.Range("M77").Value = -5507

$ws.Range("H109").Value = 49999
$ws.Range("J109").Value = 49999
$ws.Range("L109").Value = 49999
$ws.Range("N109").Value = -52773

$ws.Range("H136").Value = 2000
$ws.Range("I136").Value = 2000
$ws.Range("K136").Value = 6000
$ws.Range("M136").Value = -3450

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 96.333336
$ws.Range("I4").Value = 87.833336
$ws.Range("K4").Value = 87.833336
$ws.Range("M4").Value = 27.166664

$ws.Range("H86").Value = 2727.75
$ws.Range("I86").Value = 2727.75
$ws.Range("K86").Value = 2727.75
$ws.Range("M86").Value = -1604.75

$ws.Range("H89").Value = 2727.75
$ws.Range("I89").Value = 2727.75
$ws.Range("K89").Value = 13638.75
$ws.Range("M89").Value = -8022.75

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H6").Value = 1326
$ws.Range("I6").Value = 1439.6666
$ws.Range("J6").Value = 985
$ws.Range("K6").Value = 1439.6666
$ws.Range("L6").Value = 985
$ws.Range("M6").Value = -1326.6666
$ws.Range("N6").Value = -1211

$ws.Range("H7").Value = 66.42856999999999
$ws.Range("I7").Value = 79
$ws.Range("K7").Value = 79
$ws.Range("M7").Value = 34

$ws.Range("H17").Value = 2000
$ws.Range("J17").Value = 2850
$ws.Range("L17").Value = 2850
$ws.Range("N17").Value = -3198

$ws.Range("H25").Value = 127.75
$ws.Range("J25").Value = 100
$ws.Range("L25").Value = 100
$ws.Range("N25").Value = -448

$ws.Range("H31").Value = 10714.286
$ws.Range("I31").Value = 15000
$ws.Range("K31").Value = 15000
$ws.Range("M31").Value = -14705

$ws.Range("H34").Value = 10714.286
$ws.Range("I34").Value = 15000
$ws.Range("K34").Value = 15000
$ws.Range("M34").Value = -14798

$ws.Range("H86").Value = 12500
$ws.Range("I86").Value = 0
$ws.Range("K86").Value = 0
$ws.Range("M86").ClearContents()

$ws.Range("H89").Value = 12500
$ws.Range("I89").Value = 0
$ws.Range("K89").Value = 0
$ws.Range("M89").ClearContents()

$ws.Range("H141").Value = 874721.5
$ws.Range("J141").Value = 1146295.6
$ws.Range("L141").Value = 1146295.6
$ws.Range("N141").Value = -1156655.6

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H23").Value = 500
$ws.Range("J23").Value = 500
$ws.Range("L23").Value = 1500
$ws.Range("N23").Value = -1970

$ws.Range("H87").Value = 300
$ws.Range("I87").Value = 300
$ws.Range("K87").Value = 900
$ws.Range("M87").Value = 348

$ws.Range("H90").Value = 300
$ws.Range("I90").Value = 300
$ws.Range("K90").Value = 2700
$ws.Range("M90").Value = 3540

$ws.Range("H129").Value = 1999
$ws.Range("I129").Value = 1999
$ws.Range("J129").Value = 0
$ws.Range("K129").Value = 5997
$ws.Range("L129").Value = 0
$ws.Range("M129").Value = -997
$ws.Range("N129").ClearContents()

$ws.Range("H131").Value = 4391.3335
$ws.Range("I131").Value = 7332.8
$ws.Range("K131").Value = 21998.4
$ws.Range("M131").Value = -16958.4

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 4566.4443
$ws.Range("I46").Value = 4978.2856
$ws.Range("K46").Value = 4978.2856
$ws.Range("M46").Value = -4790.2856

$ws.Range("H82").Value = 864.8333
$ws.Range("I82").Value = 864.8333
$ws.Range("J82").Value = 0
$ws.Range("K82").Value = 864.8333
$ws.Range("L82").Value = 0
$ws.Range("M82").Value = -503.8333
$ws.Range("N82").ClearContents()

$ws.Range("H85").Value = 864.8333
$ws.Range("I85").Value = 864.8333
$ws.Range("J85").Value = 0
$ws.Range("K85").Value = 864.8333
$ws.Range("L85").Value = 0
$ws.Range("M85").Value = 383.1667
$ws.Range("N85").ClearContents()

$ws.Range("H136").Value = 2203.6
$ws.Range("I136").Value = 1704.5
$ws.Range("J136").Value = 4200
$ws.Range("K136").Value = 5113.5
$ws.Range("L136").Value = 12600
$ws.Range("M136").Value = -2563.5
$ws.Range("N136").Value = -17700

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 3095.6
$ws.Range("I113").Value = 1575
$ws.Range("J113").Value = 3475.75
$ws.Range("K113").Value = 4725
$ws.Range("L113").Value = 10427.25
$ws.Range("M113").Value = -2555
$ws.Range("N113").Value = -14767.25
